# ADD results from server
# Update row 2 values (columns A:O) across the six year sheets (2025, 2030,
# 2035, 2040, 2045, 2050 - in workbook order) with the latest results
# returned from the server.

$wb = $excel.ActiveWorkbook

$allValues = @(
    @(0, 1037.265132737054, 0, 0, 28926.05393052954, 0, 8095.925712661834, 0, 16171.06685703679, 0, 0, 48492.22142001599, 10595.37713982, 7071.74531360843, 6993.890772562212),
    @(0, 4157.588990853394, 0, 0, 45991.90904307188, 0, 8095.925712661834, 0, 37079.12819938764, 0, 0, 54844.03303316472, 17449.04999683176, 9024.733389685653, 9724.258249348202),
    @(2754.31755456332, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13034.3101291405, 12860.17168993684),
    @(2754.31755456332, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13151.8694977663, 12860.17168993684),
    @(5713.151062849596, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13601.27335063684, 14937.56878732879),
    @(5713.151062849596, 6368.910634126893, 0, 0, 57457.45307013817, 0, 8095.925712661834, 0, 52465.73681402855, 0, 0, 54844.03303316472, 21912.87293902603, 13601.27335063684, 14937.56878732879)
)

$idx = 0
foreach ($ws in $wb.Worksheets) {
    $values = $allValues[$idx]
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = $i + 1
        $ws.Cells.Item(2, $col).Value = $values[$i]
    }
    $idx = $idx + 1
}
